# Auto-generated edit script: apply value changes described by the XML diff
# (scheduled-runner price/profit recompute across multiple Leve sheets)
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 496.72726
$ws.Range("I15").Value = 496.72726
$ws.Range("K15").Value = 1490.18178
$ws.Range("M15").Value = -1321.18178
$ws.Range("H42").Value = 1108.6
$ws.Range("I42").Value = 109
$ws.Range("J42").Value = 1775
$ws.Range("K42").Value = 327
$ws.Range("L42").Value = 5325
$ws.Range("M42").Value = -97
$ws.Range("N42").Value = -5785
$ws.Range("H74").Value = 2500
$ws.Range("I74").Value = 2500
$ws.Range("K74").Value = 2500
$ws.Range("M74").Value = -1564
$ws.Range("H77").Value = 2500
$ws.Range("I77").Value = 2500
$ws.Range("K77").Value = 12500
$ws.Range("M77").Value = -7820
$ws.Range("H100").Value = 2375
$ws.Range("I100").Value = 2166.6667
$ws.Range("K100").Value = 2166.6667
$ws.Range("M100").Value = -1625.6667
$ws.Range("H135").Value = 2828.4
$ws.Range("I135").Value = 2828.4
$ws.Range("K135").Value = 25455.6
$ws.Range("M135").Value = -22920.6
$ws.Range("H137").Value = 10000
$ws.Range("I137").Value = 5000
$ws.Range("J137").Value = 15000
$ws.Range("K137").Value = 15000
$ws.Range("L137").Value = 45000
$ws.Range("M137").Value = -12450
$ws.Range("N137").Value = -50100

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5108.643
$ws.Range("I32").Value = 5108.643
$ws.Range("K32").Value = 5108.643
$ws.Range("M32").Value = -4821.643
$ws.Range("H45").Value = 5850.6665
$ws.Range("I45").Value = 4776
$ws.Range("K45").Value = 4776
$ws.Range("M45").Value = -4399
$ws.Range("H61").Value = 1708.1666
$ws.Range("I61").Value = 1708.1666
$ws.Range("K61").Value = 1708.1666
$ws.Range("M61").Value = -1496.1666
$ws.Range("H74").Value = 2557.8333
$ws.Range("I74").Value = 1471.0834
$ws.Range("K74").Value = 1471.0834
$ws.Range("M74").Value = -597.0834
$ws.Range("H77").Value = 2557.8333
$ws.Range("I77").Value = 1471.0834
$ws.Range("K77").Value = 7355.416999999999
$ws.Range("M77").Value = -2987.416999999999
$ws.Range("H97").Value = 413.6
$ws.Range("I97").Value = 342.69232
$ws.Range("K97").Value = 342.69232
$ws.Range("M97").Value = 153.30768
$ws.Range("H113").Value = 398
$ws.Range("J113").Value = 398
$ws.Range("L113").Value = 398
$ws.Range("N113").Value = -9076
$ws.Range("H118").Value = 0
$ws.Range("J118").Value = 0
$ws.Range("L118").Value = 0
$ws.Range("N118").ClearContents()
$ws.Range("H122").Value = 1743.6666
$ws.Range("I122").Value = 1743.6666
$ws.Range("K122").Value = 5230.9998
$ws.Range("M122").Value = -2780.9998
$ws.Range("H132").Value = 4359.75
$ws.Range("I132").Value = 2989.625
$ws.Range("K132").Value = 8968.875
$ws.Range("M132").Value = -6438.875
$ws.Range("H136").Value = 1708.1666
$ws.Range("I136").Value = 1708.1666
$ws.Range("K136").Value = 5124.4998
$ws.Range("M136").Value = -2574.4998

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 3926.842
$ws.Range("I86").Value = 3779.6428
$ws.Range("K86").Value = 3779.6428
$ws.Range("M86").Value = -2656.6428
$ws.Range("H89").Value = 3926.842
$ws.Range("I89").Value = 3779.6428
$ws.Range("K89").Value = 18898.214
$ws.Range("M89").Value = -13282.214
$ws.Range("H99").Value = 1530.1666
$ws.Range("I99").Value = 1514
$ws.Range("K99").Value = 1514
$ws.Range("M99").Value = -16
$ws.Range("H134").Value = 2659.5
$ws.Range("I134").Value = 2591.4
$ws.Range("K134").Value = 7774.200000000001
$ws.Range("M134").Value = -5239.200000000001
$ws.Range("H141").Value = 0
$ws.Range("I141").Value = 0
$ws.Range("K141").Value = 0
$ws.Range("M141").ClearContents()

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H60").Value = 18749.5
$ws.Range("J60").Value = 0
$ws.Range("L60").Value = 0
$ws.Range("N60").ClearContents()
$ws.Range("H122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("M122").ClearContents()
$ws.Range("H132").Value = 3428.625
$ws.Range("I132").Value = 3737.8
$ws.Range("K132").Value = 11213.4
$ws.Range("M132").Value = -8683.400000000001
$ws.Range("H134").Value = 3992.3333
$ws.Range("I134").Value = 3992.3333
$ws.Range("K134").Value = 11976.9999
$ws.Range("M134").Value = -9441.999899999999

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H138").Value = 1900
$ws.Range("J138").Value = 2000
$ws.Range("L138").Value = 6000
$ws.Range("N138").Value = -16280

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 107.125
$ws.Range("J2").Value = 134.66667
$ws.Range("L2").Value = 134.66667
$ws.Range("N2").Value = -360.66667
$ws.Range("H11").Value = 11677647
$ws.Range("I11").Value = 15863637
$ws.Range("K11").Value = 15863637
$ws.Range("M11").Value = -15863498
$ws.Range("H14").Value = 78065.89999999999
$ws.Range("I14").Value = 250150
$ws.Range("J14").Value = 35044.875
$ws.Range("K14").Value = 250150
$ws.Range("L14").Value = 35044.875
$ws.Range("M14").Value = -249982
$ws.Range("N14").Value = -35380.875
$ws.Range("H80").Value = 13426.286
$ws.Range("I80").Value = 3047.5
$ws.Range("J80").Value = 17577.8
$ws.Range("K80").Value = 3047.5
$ws.Range("L80").Value = 17577.8
$ws.Range("M80").Value = -2049.5
$ws.Range("N80").Value = -19573.8
$ws.Range("H83").Value = 13426.286
$ws.Range("I83").Value = 3047.5
$ws.Range("J83").Value = 17577.8
$ws.Range("K83").Value = 15237.5
$ws.Range("L83").Value = 87889
$ws.Range("M83").Value = -10245.5
$ws.Range("N83").Value = -97873
$ws.Range("H97").Value = 977.1429000000001
$ws.Range("I97").Value = 890
$ws.Range("K97").Value = 890
$ws.Range("M97").Value = -394
$ws.Range("H122").Value = 1585.091
$ws.Range("I122").Value = 1715.1111
$ws.Range("J122").Value = 1000
$ws.Range("K122").Value = 5145.3333
$ws.Range("L122").Value = 3000
$ws.Range("M122").Value = -2695.3333
$ws.Range("N122").Value = -7900

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H17").Value = 6000
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 6000
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 6000
$ws.Range("M17").ClearContents()
$ws.Range("N17").Value = -6340
$ws.Range("H46").Value = 402.8889
$ws.Range("I46").Value = 416.2857
$ws.Range("K46").Value = 416.2857
$ws.Range("M46").Value = -228.2857
$ws.Range("H93").Value = 714.5714
$ws.Range("I93").Value = 714.5714
$ws.Range("K93").Value = 714.5714
$ws.Range("M93").Value = 533.4286
$ws.Range("H122").Value = 3764
$ws.Range("I122").Value = 3764
$ws.Range("K122").Value = 11292
$ws.Range("M122").Value = -8842

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H29").Value = 0
$ws.Range("I29").Value = 0
$ws.Range("J29").Value = 0
$ws.Range("K29").Value = 0
$ws.Range("L29").Value = 0
$ws.Range("M29").ClearContents()
$ws.Range("N29").ClearContents()
$ws.Range("H122").Value = 696.75
$ws.Range("I122").Value = 625
$ws.Range("K122").Value = 1875
$ws.Range("M122").Value = 575
$ws.Range("H136").Value = 9950.25
$ws.Range("I136").Value = 9950.25
$ws.Range("K136").Value = 29850.75
$ws.Range("M136").Value = -27300.75
